$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Remove the "Laboratory Report Narrative" (row 18) and
# "Pathology Report Narrative" (row 19) rows from the Clinical Notes
# section. Deleting both rows in one call shifts everything below up by two.
$ws.Range("A18:J19").EntireRow.Delete()
